# Applies the update described by the commit "Upload new version with timestamp":
#  - Inserts a new product "فرش اسنان فوكس " into the price list. It takes row 94,
#    right after "فرش اسنان اورل فريش", and every later product row's data shifts
#    down by one slot (with a couple of entries also getting refreshed figures).
#  - Updates the running totals row and the footer row accordingly (they move from
#    rows 103/104 to rows 104/105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row at 103. This pushes the old totals row (103 -> 104) and
#    the old footer row (104 -> 105) down, together with their merged cells.
#    Then copy the formatting of row 102 onto the freshly inserted row 103 so
#    it matches the other product rows (same style ids as rows 94-102).
# ---------------------------------------------------------------------------
$ws.Rows(103).Insert()

$ws.Range("A102:N102").Copy()
$ws.Range("A103:N103").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows(103).RowHeight = 25.5

# ---------------------------------------------------------------------------
# 2) Write the final contents for the product rows 94..103.
#    Columns: A = running number, B = product name, H = "current balance",
#    L = sale price, N = number of transactions.
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row = 94;  A = 91;  B = "فرش اسنان فوكس ";          H = "1:0";  L = 60; N = "1:0" },
    @{ Row = 95;  A = 92;  B = "كالونا ";                     H = "-1:0"; L = 15; N = "1:0" },
    @{ Row = 96;  A = 93;  B = "كريم ONE للبشره الحساسه";      H = "19:0"; L = 50; N = "2:0" },
    @{ Row = 97;  A = 94;  B = "كريم براشوت";                  H = "0:0";  L = 45; N = "1:0" },
    @{ Row = 98;  A = 95;  B = "ماء اكسجين 20";                H = "0:0";  L = 10; N = "1:0" },
    @{ Row = 99;  A = 96;  B = "ماسك جلسات اطفال";             H = "-1:0"; L = 20; N = "1:0" },
    @{ Row = 100; A = 97;  B = "معجون سيجنال 120 مل ";         H = "5:0";  L = 60; N = "1:0" },
    @{ Row = 101; A = 98;  B = "معجون سيجنال 25 مل";           H = "5:0";  L = 40; N = "2:0" },
    @{ Row = 102; A = 99;  B = "معجون سيجنال عرض 50ملل";       H = "3:0";  L = 80; N = "2:0" },
    @{ Row = 103; A = 100; B = "معجون كلوز اب الصغير";         H = "16:0"; L = 20; N = "1:0" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value2 = $item.A
    $ws.Cells.Item($r, 2).Value2 = $item.B
    $ws.Cells.Item($r, 8).Value2 = $item.H
    $ws.Cells.Item($r, 12).Value2 = $item.L
    $ws.Cells.Item($r, 14).Value2 = $item.N
}

# ---------------------------------------------------------------------------
# 3) Update the running-total row (now row 104) and the footer row (now row
#    105, whose height shrinks slightly to match the refreshed export).
# ---------------------------------------------------------------------------
$ws.Cells.Item(104, 11).Value2 = 4807.5200000000004

$ws.Rows(105).RowHeight = 16.5
